$d = $word.ActiveDocument

# 1. Update activation date
$d.Content.Find.Execute(
    "Ativação: 01/01/2018", $true, $false, $false, $false, $false,
    $true, 1, $false, "Ativação: 01/01/2025", 2
)

# 2. Update Portuguese short syllabus ("Programa resumido")
$d.Content.Find.Execute(
    "Espaços vetoriais, Transformações lineares, auto-valores e auto-vetores, diagonalização de Operadores, espaços vetoriais com produto interno, aplicações as equações diferenciais.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Espaços vetoriais, transformações lineares, diagonalização de operadores lineares, forma canônica de Jordan, espaços vetoriais com produto interno, aplicações a sistemas dinâmicos.",
    2
)

# 3. Update English short syllabus (italic)
$d.Content.Find.Execute(
    "Vector Spaces, linear Transformations, eigenvalues and eigenvectors, Diagonalization, Inner product in vectorial Spaces, applications to differential Equations.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Vector spaces, linear transformations, diagonalization of linear operators, Jordan canonical form, inner product spaces, applications to dynamical systems.",
    2
)

# 4. Fix typo "spaços" -> "Espaços" at the start of the full "Programa" paragraph
$d.Content.Find.Execute(
    "spaços vetoriais: Definição, Propriedades dos Espaços Vetoriais",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Espaços vetoriais: Definição, Propriedades dos Espaços Vetoriais",
    2
)

# 5. Update Bibliografia list.
# This run's text is 566 chars long, well beyond Word's ~255-char Find/Replace
# text limits, so locate the paragraph directly (the one right after the
# "Bibliografia" Heading2) and replace its text via Range instead.
$bibHeadingIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq "Bibliografia") {
        $bibHeadingIndex = $i
        break
    }
}

if ($bibHeadingIndex -ge 1) {
    $bibPara = $d.Paragraphs.Item($bibHeadingIndex + 1)
    $bibRange = $bibPara.Range
    # Exclude the trailing paragraph mark so we only replace the run text.
    $bibRange.MoveEnd(1, -1)
    $bibRange.Text = "1.LIMA, Elon Lages. Álgebra Linear, IMPA, 2020. ISBN: 978-65-990528-3-5. 10ª edição.2.HOFFMAN, Kennethe; KUNZE, Ray. Linear Algebra. Pearson. 1971. 2nd Edition.3. STRANG, Gilbert. Álgebra linear e suas aplicações, São Paulo: Cengage Learning, 2010.4.LIPSCHUTZ, Seymour. Álgebra linear. 3. ed. São Paulo: Ed. McGrawHill. 1990.5.HOWARD, Anton ; RORRES, Chris. Álgebra linear com aplicações. 8. ed., Ed. Bookman, 2001.6.MICHOLSON, W. Keith. Álgebra linear. 2. ed. São Paulo: Ed. Mc GrawHill, 2006.7.BOLDRINI, José Luiz ; COSTA Sueli I. Rodrigues; FIGUEIREDO Vera Lúcia; WETZLER Henry G. Álgebra linear., 3. ed. São Paulo: Editora Harbra Ltda, 1986.8.POOLE, David. Álgebra linear. São Paulo: Pioneira Thomson Learning, 2004."
}
